# Insert a new weekly price record as row 29 (Berenjena, Terminal La Palmera
# de La Serena). All subsequent rows (old 29..120) shift down by one row to
# 30..121, preserving their original values, and the sheet dimension grows
# from A1:R120 to A1:R121.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing row 29 (and everything below it) down by one row.
$ws.Rows.Item(29).Insert()

# Populate the newly inserted row 29 with the new record's data.
$ws.Range("A29").Value = 8
$ws.Range("B29").Value = "Terminal La Palmera de La Serena"
$ws.Range("C29").Value = "Coquimbo"
$ws.Range("D29").Value = 44620
$ws.Range("E29").Value = 4
$ws.Range("F29").Value = 100112001
$ws.Range("G29").Value = "Berenjena"
$ws.Range("H29").Value = "Sin especificar"
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 520
$ws.Range("K29").Value = 8000
$ws.Range("L29").Value = 9000
$ws.Range("M29").Value = 8500
$ws.Range("N29").Value = "`$/caja 50 unidades"
$ws.Range("O29").Value = "Región de Arica y Parinacota"
$ws.Range("P29").Value = 170
$ws.Range("Q29").Value = 50
$ws.Range("R29").Value = "Hortaliza"
